$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.253.68"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "'1.827.09"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  -0.95%  "
$ws.Range("D5").Value = "'313.71"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("E7").Value = "  -1.89%  "
$ws.Range("D8").Value = "'0.3710"
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").Value = "'0.07257"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").Value = "'0.8648"
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("D12").Value = "'1.832.58"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").Value = "'6.738"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").Value = "'5.325"
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("D15").Value = "'0.07090"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "'89.50"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").Value = "'1.005"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("D18").Value = "'0.000008875"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("D21").Value = "'27.361.06"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").Value = "'5.145"
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("D24").Value = "'2.057.36"
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("D25").Value = "'1.993"
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("D26").Value = "'152.58"
$ws.Range("E26").Value = "  -2.00%  "
$ws.Range("D27").Value = "'2.209"
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("D28").Value = "'18.43"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").Value = "'5.249"
$ws.Range("E29").Value = "  -3.06%  "
$ws.Range("D30").Value = "'116.67"
$ws.Range("E30").Value = "  -3.15%  "
$ws.Range("D31").Value = "'0.08865"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").Value = "'1.200"
$ws.Range("E32").Value = "  -2.89%  "
$ws.Range("D33").Value = "'0.7603"
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("D34").Value = "'4.464"
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("D35").Value = "'2.802"
$ws.Range("E35").Value = "  -4.12%  "
$ws.Range("D36").Value = "'1.004"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").Value = "'1.121"
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("D38").Value = "'0.01978"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").Value = "'0.05268"
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("D40").Value = "'7.342"
$ws.Range("E40").Value = "  +2.62%  "
$ws.Range("D41").Value = "'2.871"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "'0.1701"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("D44").Value = "'8.731"
$ws.Range("E44").Value = "  -2.40%  "
$ws.Range("D45").Value = "'10.62"
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("D46").Value = "'107.64"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").Value = "'0.4761"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").Value = "'0.06396"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("D50").Value = "'1.676"
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("D51").Value = "'1.868"
$ws.Range("E51").Value = "  -1.53%  "
